# Update version/build strings for release "mines - version 1.0.0 (Feb 3 2026)"
# Old build stamp: "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# New build stamp: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Jharia Coal Mine, India, M1686, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

$wsData.Range("S2").Value = $newVersion
$wsData.Range("S3").Value = $newVersion
$wsData.Range("S4").Value = $newVersion
$wsData.Range("S5").Value = $newVersion
$wsData.Range("S6").Value = $newVersion
$wsData.Range("S7").Value = $newVersion
